$wb = $excel.ActiveWorkbook

# New row data for each sheet (row 69), derived from row 68 pattern of next-day log entry.
$rows = @{
    "MID_LFT_#1" = @{
        A = 45855.4646412037
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
        D = "0x01,0x4C"
        E = "0x07"
        F = 400
        G = [double]"5.68631262647113e+23"
        H = 332
        I = 7
    }
    "MID_LFT_#2" = @{
        A = 45855.4646412037
        B = "0x01,0x7c"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
        D = "0x01,0x4C"
        E = "0x19"
        F = 380
        G = [double]"5.68432987514711e+23"
        H = 332
        I = 25
    }
    "MID_PLT_#1" = @{
        A = 45855.4646412037
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
        D = "0x00,0x65"
        E = "0x15"
        F = 110
        G = [double]"5.68631262647113e+23"
        H = 101
        I = 15
    }
    "MID_PLT_#2" = @{
        A = 45855.4646412037
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
        D = "0x00,0x79"
        E = "0x9"
        F = 130
        G = [double]"5.68631262647113e+23"
        H = 121
        I = 9
    }
}

foreach ($sheetName in $rows.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $data = $rows[$sheetName]
    $r = 69

    $ws.Cells.Item($r, 1).Value = $data.A
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 2).Value = $data.B
    $ws.Cells.Item($r, 3).Value = $data.C
    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 5).Value = $data.E

    $ws.Cells.Item($r, 6).Value = $data.F
    $ws.Cells.Item($r, 7).Value = $data.G
    $ws.Cells.Item($r, 8).Value = $data.H
    $ws.Cells.Item($r, 9).Value = $data.I
}
